$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Context: the "ReportManagerController.cs" section of the readme.
# Before:
#   P28: "If you are only using ... endpoints." + " A new compiler symbol ... MOBILIZER_ENDPOINT. "
#   P29: "If you comment out the line:"
#   P30: "//#define MOBILIZER_ENDPOINT"
#   P31: (empty)
#   P32: "You will get default implementations ... as " + bookmark(_GoBack) + "easy as possible."
#   P33: (empty, trailing)
#
# After:
#   P28: "In some rare cases ... Report Server DB."                    (new paragraph)
#   P29: "So assuming you are only using ... GetMobilizerSettings endpoint. "  (new paragraph)
#   P30: "A new compiler symbol ... MOBILIZER_ENDPOINT. "              (new paragraph)
#   P31: (empty, now holds bookmark(_GoBack))                          (new paragraph)
#   P32: "If you comment out the line:"                                (shifted, unchanged)
#   P33: "//#define MOBILIZER_ENDPOINT"                                (shifted, unchanged)
#   P34: (empty)                                                       (shifted, unchanged)
#   P35: "You will get generic (do nothing) implementations ... easy as possible." (bookmark removed)
#   (trailing empty paragraph removed)
# ---------------------------------------------------------------------------

# --- Step 1: rewrite paragraph 28's text in place (single literal Find/Replace
#     over the whole paragraph keeps the "highlight=white" run formatting). ---
$p28 = $d.Paragraphs(28)
$oldP28 = "If you are only using the Forerunner " + [char]0x201c + "reportViewerEZ" + [char]0x201c + "and not the " + [char]0x201c + "reportExplorerEZ" + [char]0x201d + " widget. You will still need to provide a default implementation for the report manager endpoints. A new compiler symbol has been defined at the top of ReportManagerController.cs named MOBILIZER_ENDPOINT. "
$newP28 = "In some rare cases you may not want to use the standard Mobilizer endpoints defined in the Report Manager Controller. This may be the case if you don" + [char]0x2019 + "t want forerunner SDK to read / write forerunner specific data into the Report Server DB."
$found = $p28.Range.Find.Execute($oldP28, $false, $false, $false, $false, $false, $true, 1, $false, $newP28, 2)
if (-not $found) {
    throw "Step 1: could not find/replace paragraph 28 text"
}

# --- Step 2: insert three new empty paragraphs right after paragraph 28,
#     matching the "empty paragraph w/ highlight=white pilcrow" styling. ---
$p28.Range.InsertParagraphAfter()
$p28.Range.InsertParagraphAfter()
$p28.Range.InsertParagraphAfter()

# --- Step 3: fill in the 3 new paragraphs (29, 30, 31). ---
$newP29 = "So assuming you are only using the Forerunner " + [char]0x201c + "reportViewerEZ" + [char]0x201c + "and not the " + [char]0x201c + "reportExplorerEZ" + [char]0x201d + " widget. You will still need to provide an implementation for the report manager endpoints because some of these endpoints are used by " + [char]0x201c + "reportViewerEZ" + [char]0x201d + ". One example of this it the GetMobilizerSettings endpoint. "
$p29 = $d.Paragraphs(29)
$p29.Range.InsertAfter($newP29)

$newP30 = "A new compiler symbol has been defined at the top of ReportManagerController.cs named MOBILIZER_ENDPOINT. "
$p30 = $d.Paragraphs(30)
$p30.Range.InsertAfter($newP30)

# Paragraph 31 stays empty for now -- the _GoBack bookmark is added to it in Step 5.

# --- Step 4: rewrite the final content paragraph (now shifted to 35); this
#     also removes the old bookmark that used to sit inside it. ---
$pLast = $d.Paragraphs(35)
$oldLast = "You will get default implementations for all the report manager endpoints. You may want to implement select endpoints yourself such as get and save thumbnail. It is recommended you put your implementations into this file and carefully comment your changes so that future upgrades are made as easy as possible."
$newLast = "You will get generic (do nothing) implementations for all the report manager endpoints. You may want to implement select endpoints yourself such as get and save thumbnail. It is recommended you put your implementations into this file and carefully comment your changes so that future upgrades are made as easy as possible."
$foundLast = $pLast.Range.Find.Execute("You will get*easy as possible.", $false, $false, $true, $false, $false, $true, 1, $false, $newLast, 2)
if (-not $foundLast) {
    throw "Step 4: could not find/replace final paragraph text"
}

# --- Step 5: (re)place the _GoBack bookmark on the new empty paragraph 31.
#     Word only keeps a single _GoBack bookmark, so adding it here removes
#     any other instance automatically. ---
$p31 = $d.Paragraphs(31)
$bmRange = $p31.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 6: delete the trailing empty paragraph before the section break
#     (now shifted to paragraph 36). ---
$trailing = $d.Paragraphs(36)
$trailing.Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
